# Insert a brand-new data row at row 46 (pushing the existing rows 46-127
# down to 47-128), then populate the new row with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(46).Insert()

$ws.Range("A46").Value = 9
$ws.Range("B46").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C46").Value = 'Metropolitana'
$ws.Range("D46").Value = 44791
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = 100112022
$ws.Range("G46").Value = 'Arveja Verde'
$ws.Range("H46").Value = 'Sin especificar'
$ws.Range("I46").Value = 'Primera'
$ws.Range("J46").Value = 34
$ws.Range("K46").Value = 42000
$ws.Range("L46").Value = 42000
$ws.Range("M46").Value = 42000
$ws.Range("N46").Value = '$/malla 25 kilos'
$ws.Range("O46").Value = 'Provincia de Limarí'
$ws.Range("P46").Value = 1680
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = 'Hortaliza'
